$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.02492466666666666
$ws.Range("H2").Value = 0.07477399999999999
$ws.Range("I2").Value = 0.5549276228148621
$ws.Range("J2").Value = 0.6515968803102261
$ws.Range("M2").Value = 35.418552
$ws.Range("N2").Value = 70.837104
$ws.Range("O2").Value = 0.1918645458500658
$ws.Range("P2").Value = 0.1411406244578309
$ws.Range("Q2").Value = 0.8827956024159999
$ws.Range("R2").Value = 5.296773614495999
$ws.Range("S2").Value = 0.1064709363310301
$ws.Range("T2").Value = 0.09196679058175983
$ws.Range("G3").Value = 0.02492466666666666
$ws.Range("H3").Value = 0.07477399999999999
$ws.Range("I3").Value = 0.5549276228148621
$ws.Range("J3").Value = 0.6515968803102261
$ws.Range("O3").Value = 0.00307846483153751
$ws.Range("P3").Value = 0.003396899985635307
$ws.Range("S3").Value = 0.001708325170884265
$ws.Range("T3").Value = 0.002213409433365818
$ws.Range("G4").Value = 0.02492466666666666
$ws.Range("H4").Value = 0.07477399999999999
$ws.Range("I4").Value = 0.5549276228148621
$ws.Range("J4").Value = 0.6515968803102261
$ws.Range("M4").Value = 65.80493666666666
$ws.Range("N4").Value = 197.41481
$ws.Range("O4").Value = 0.356469521629296
$ws.Range("P4").Value = 0.3933425844261512
$ws.Range("Q4").Value = 1.640166111437777
$ws.Range("R4").Value = 14.76149500294
$ws.Range("S4").Value = 0.1978147842436963
$ws.Range("T4").Value = 0.2563008009052418
$ws.Range("G5").Value = 0.02492466666666666
$ws.Range("H5").Value = 0.07477399999999999
$ws.Range("I5").Value = 0.5549276228148621
$ws.Range("J5").Value = 0.6515968803102261
$ws.Range("M5").Value = 16.4967725
$ws.Range("N5").Value = 32.993545
$ws.Range("O5").Value = 0.08936406445143084
$ws.Range("P5").Value = 0.06573856469877065
$ws.Range("Q5").Value = 0.4111765556383333
$ws.Range("R5").Value = 2.46705933383
$ws.Range("S5").Value = 0.04959058785110664
$ws.Range("T5").Value = 0.04283504367379092
$ws.Range("G6").Value = 0.02492466666666666
$ws.Range("H6").Value = 0.07477399999999999
$ws.Range("I6").Value = 0.5549276228148621
$ws.Range("J6").Value = 0.6515968803102261
$ws.Range("M6").Value = 27.40387933333333
$ws.Range("N6").Value = 82.21163799999999
$ws.Range("O6").Value = 0.1484485549499597
$ws.Range("P6").Value = 0.1638040132897181
$ws.Range("Q6").Value = 0.6830325577568888
$ws.Range("R6").Value = 6.147293019811999
$ws.Range("S6").Value = 0.08237820370868255
$ws.Range("T6").Value = 0.1067341840418751
$ws.Range("G7").Value = 0.02492466666666666
$ws.Range("H7").Value = 0.07477399999999999
$ws.Range("I7").Value = 0.5549276228148621
$ws.Range("J7").Value = 0.6515968803102261
$ws.Range("M7").Value = 38.909429
$ws.Range("N7").Value = 116.728287
$ws.Range("O7").Value = 0.2107748482877103
$ws.Range("P7").Value = 0.232577313141894
$ws.Range("Q7").Value = 0.9698045480153331
$ws.Range("R7").Value = 8.728240932137998
$ws.Range("S7").Value = 0.1169647855094623
$ws.Range("T7").Value = 0.1515466516741927
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.0199905
$ws.Range("H8").Value = 0.039981
$ws.Range("I8").Value = 0.4450723771851379
$ws.Range("J8").Value = 0.3484031196897739
$ws.Range("M8").Value = 35.418552
$ws.Range("N8").Value = 70.837104
$ws.Range("O8").Value = 0.1918645458500658
$ws.Range("P8").Value = 0.1411406244578309
$ws.Range("Q8").Value = 0.708034563756
$ws.Range("R8").Value = 2.832138255024
$ws.Range("S8").Value = 0.08539360951903568
$ws.Range("T8").Value = 0.0491738338760711
$ws.Range("E9").Value = 2
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.0199905
$ws.Range("H9").Value = 0.039981
$ws.Range("I9").Value = 0.4450723771851379
$ws.Range("J9").Value = 0.3484031196897739
$ws.Range("O9").Value = 0.00307846483153751
$ws.Range("P9").Value = 0.003396899985635307
$ws.Range("Q9").Value = 0.0113604079085
$ws.Range("R9").Value = 0.068162447451
$ws.Range("S9").Value = 0.001370139660653245
$ws.Range("T9").Value = 0.001183490552269489
$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.0199905
$ws.Range("H10").Value = 0.039981
$ws.Range("I10").Value = 0.4450723771851379
$ws.Range("J10").Value = 0.3484031196897739
$ws.Range("M10").Value = 65.80493666666666
$ws.Range("N10").Value = 197.41481
$ws.Range("O10").Value = 0.356469521629296
$ws.Range("P10").Value = 0.3933425844261512
$ws.Range("Q10").Value = 1.315473586435
$ws.Range("R10").Value = 7.89284151861
$ws.Range("S10").Value = 0.1586547373855997
$ws.Range("T10").Value = 0.1370417835209093
$ws.Range("E11").Value = 2
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 0.0199905
$ws.Range("H11").Value = 0.039981
$ws.Range("I11").Value = 0.4450723771851379
$ws.Range("J11").Value = 0.3484031196897739
$ws.Range("M11").Value = 16.4967725
$ws.Range("N11").Value = 32.993545
$ws.Range("O11").Value = 0.08936406445143084
$ws.Range("P11").Value = 0.06573856469877065
$ws.Range("Q11").Value = 0.32977873066125
$ws.Range("R11").Value = 1.319114922645
$ws.Range("S11").Value = 0.0397734766003242
$ws.Range("T11").Value = 0.02290352102497974
$ws.Range("E12").Value = 2
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 0.0199905
$ws.Range("H12").Value = 0.039981
$ws.Range("I12").Value = 0.4450723771851379
$ws.Range("J12").Value = 0.3484031196897739
$ws.Range("M12").Value = 27.40387933333333
$ws.Range("N12").Value = 82.21163799999999
$ws.Range("O12").Value = 0.1484485549499597
$ws.Range("P12").Value = 0.1638040132897181
$ws.Range("Q12").Value = 0.5478172498130001
$ws.Range("R12").Value = 3.286903498878
$ws.Range("S12").Value = 0.06607035124127712
$ws.Range("T12").Value = 0.05706982924784297
$ws.Range("E13").Value = 2
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 0.0199905
$ws.Range("H13").Value = 0.039981
$ws.Range("I13").Value = 0.4450723771851379
$ws.Range("J13").Value = 0.3484031196897739
$ws.Range("M13").Value = 38.909429
$ws.Range("N13").Value = 116.728287
$ws.Range("O13").Value = 0.2107748482877103
$ws.Range("P13").Value = 0.232577313141894
$ws.Range("Q13").Value = 0.7778189404244999
$ws.Range("R13").Value = 4.666913642547
$ws.Range("S13").Value = 0.09381006277824804
$ws.Range("T13").Value = 0.08103066146770131
